# Rename the "RMM" index row (row 14 of the Plan1 table) to "MJO".
#   Index      (col A): "RMM "                -> "MJO"
#   Index_Unit (col C): "RMM (dimensionless)"  -> "MJO (dimensionless)"
# The Unit, Name_Index, Methodology, Access and Reference columns for that
# row are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "MJO"
$ws.Range("C14").Value = "MJO (dimensionless)"

# Reflect the author's final selection/view in the saved file.
$ws.Range("D14").Select()
